$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "User Needs Chart" paragraph right after the title
#    paragraph and before the table.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newParaXml = '<w:p><w:pPr><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">User Needs Chart </w:t></w:r></w:p>'
$newPara.Range.InsertXML($newParaXml)

# ---------------------------------------------------------------------------
# 2) "Mass add students" -> "Mass " + proofErr(gramStart) + "add" +
#    proofErr(gramEnd) + " students"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Mass add students") | Out-Null
$xml = '<w:p w14:paraId="5E3E7296" w14:textId="246B9ED7" w:rsidR="000B0736" w:rsidRPr="000B0736" w:rsidRDefault="00D32DA7" w:rsidP="000B0736">' +
       '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Mass </w:t></w:r>' +
       '<w:proofErr w:type="gramStart"/>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>add</w:t></w:r>' +
       '<w:proofErr w:type="gramEnd"/>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> students</w:t></w:r>' +
       '</w:p>'
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 3) "Mass add professors" -> "Mass " + proofErr(gramStart) + "add" +
#    proofErr(gramEnd) + " professors"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Mass add professors") | Out-Null
$xml = '<w:p w14:paraId="544E9B51" w14:textId="59AFFF64" w:rsidR="000B0736" w:rsidRPr="000B0736" w:rsidRDefault="00D32DA7" w:rsidP="000B0736">' +
       '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">Mass </w:t></w:r>' +
       '<w:proofErr w:type="gramStart"/>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>add</w:t></w:r>' +
       '<w:proofErr w:type="gramEnd"/>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> professors</w:t></w:r>' +
       '</w:p>'
$rng.InsertXML($xml)

# ---------------------------------------------------------------------------
# 4) Add explicit "no shading" <w:shd .../> to the "Take exams" and
#    "Upload assignments" cells (row 7, columns 2 & 3 of the one table).
# ---------------------------------------------------------------------------
$t = $d.Tables(1)
foreach ($colIdx in 2, 3) {
    $cell = $t.Cell(7, $colIdx)
    $cell.Shading.Texture = 0
    $cell.Shading.ForegroundPatternColor = -16777216
    $cell.Shading.BackgroundPatternColor = -16777216
}

# ---------------------------------------------------------------------------
# 5) "Access to cheating tools" -> "Access to cheating" + " prevention" +
#    " tools" (three runs, identical formatting).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Access to cheating tools") | Out-Null
$xml = '<w:p w14:paraId="6D286B2C" w14:textId="5CBE6202" w:rsidR="00D32DA7" w:rsidRDefault="00D32DA7" w:rsidP="000B0736">' +
       '<w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Access to cheating</w:t></w:r>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> prevention</w:t></w:r>' +
       '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> tools</w:t></w:r>' +
       '</w:p>'
$rng.InsertXML($xml)

Write-Output "edit.ps1 completed"
